$wb = $excel.ActiveWorkbook

$wsYearly = $wb.Worksheets.Item("Yearly")
$wsAllTime = $wb.Worksheets.Item("All Time")

# Update the December Taxable Account dividend value for 2016 (Yearly sheet)
$wsYearly.Range("D14").Value = 133.62

# Update selections to match final workbook state
$wsYearly.Activate()
$wsYearly.Range("M9").Select()

$wsAllTime.Activate()
$wsAllTime.Range("P26").Select()
